# Applies the "add cap statements for all actors" edit to alert-sender.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix typo + rewrite security note on the "meta" sheet
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Range("B3").Value = "This Section describes the expected capabilities of the Da Vinci Alert Sender actor which is responsible for sending the alert, typically operated by the facility or organization where the event occurred. The complete list of FHIR profiles, RESTful operations, and search parameters supported by Da Vinci Alert Senders are defined."
$wsMeta.Range("B7").Value = "1. For general security consideration refer to the [Security and Privacy Considerations](http://build.fhir.org/secpriv-module.html). `n1. For security considerations specific to this guide refer to the  [Security](security.html) page for requirements and recommendations.`n1. A server **SHALL** reject any unauthorized requests by returning an `HTTP 401` unauthorized response code."

# ---------------------------------------------------------------------------
# 2. Restructure the "interactions" sheet: drop the per-resource columns for
#    AllergyIntolerance / CarePlan / CareTeam / Condition (conformance is now
#    tracked elsewhere), leaving Device onward.
# ---------------------------------------------------------------------------
$wsInteractions = $wb.Worksheets.Item("interactions")
$wsInteractions.Range("B1:E1").EntireColumn.Delete()
$wsInteractions.Range("C13").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "rest_interactions" sheet (system/type level REST
#    interactions: transaction, batch, search-system, history-system) between
#    "interactions" and "sps".
# ---------------------------------------------------------------------------
$wsSps = $wb.Worksheets.Item("sps")
$wsRest = $wb.Worksheets.Add($wsSps)
$wsRest.Name = "rest_interactions"

$wsRest.Range("A1").Value = "code"
$wsRest.Range("B1").Value = "conf"
$wsRest.Range("C1").Value = "doc"

$wsRest.Range("A2").Value = "transaction"
$wsRest.Range("B2").Value = "SHALL"
$wsRest.Range("C2").Value = "Whether as a direct push based transaction or via subscription notification, a common  ``transaction``  type " + [char]0x201C + "Alert Bundle" + [char]0x201D + " is the FHIR object that is exchanged between the Da Vinci Alert Actors."

$wsRest.Range("A3").Value = "batch"
$wsRest.Range("B3").Value = "MAY"

$wsRest.Range("A4").Value = "search-system"
$wsRest.Range("B4").Value = "MAY"

$wsRest.Range("A5").Value = "history-system"
$wsRest.Range("B5").Value = "MAY"

$wsRest.Rows.Item(2).RowHeight = 90
$wsRest.Range("G1").WrapText = $true
$wsRest.Range("G2").WrapText = $true

$wsRest.Range("C2").Select()

# ---------------------------------------------------------------------------
# 4. Restore cosmetic selections on the other sheets and make "meta" the
#    active tab (it was "igs" before).
# ---------------------------------------------------------------------------
$wsIgs = $wb.Worksheets.Item("igs")
$wsIgs.Range("C4").Select()

$wsOps = $wb.Worksheets.Item("ops")
$wsOps.Range("C4").Select()

$wsMeta.Activate()
$wsMeta.Range("B12").Select()

Write-Output "done"
